$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.803.05"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.750.41"
$ws.Range("E3").Value = "  -4.68%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.87"
$ws.Range("E5").Value = "  -7.71%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5073"
$ws.Range("E7").Value = "  -4.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.67"
$ws.Range("E8").Value = "  -6.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2640"
$ws.Range("E9").Value = "  -6.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06168"
$ws.Range("E10").Value = "  -11.18%  "
$ws.Range("D11").Value = "1.756.13"
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.69"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.06893"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6031"
$ws.Range("E14").Value = "  -14.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.493"
$ws.Range("E15").Value = "  -8.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.05"
$ws.Range("E16").Value = "  -11.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "25.812.47"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  -11.62%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006840"
$ws.Range("E21").Value = "  -6.79%  "
$ws.Range("D22").Value = "1.980.64"
$ws.Range("E22").Value = "  -5.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.082"
$ws.Range("E23").Value = "  -9.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.217"
$ws.Range("E24").Value = "  -8.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.195"
$ws.Range("E25").Value = "  -11.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.49"
$ws.Range("E26").Value = "  -3.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.471"
$ws.Range("E27").Value = "  -12.25%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.817"
$ws.Range("E28").Value = "  -11.57%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.02"
$ws.Range("E29").Value = "  -9.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.98"
$ws.Range("E30").Value = "  -5.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08185"
$ws.Range("E31").Value = "  -6.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.708"
$ws.Range("E32").Value = "  -9.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.483"
$ws.Range("E33").Value = "  -10.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04502"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9998"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.647"
$ws.Range("E36").Value = "  -8.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9939"
$ws.Range("E37").Value = "  -10.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6061"
$ws.Range("E38").Value = "  -13.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.699"
$ws.Range("E39").Value = "  -12.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.922"
$ws.Range("E41").Value = "  -12.35%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.33"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3827"
$ws.Range("E44").Value = "  -15.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7396"
$ws.Range("E45").Value = "  -15.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.945"
$ws.Range("E46").Value = "  -14.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05449"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1105"
$ws.Range("E48").Value = "  -7.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.956"
$ws.Range("E49").Value = "  -15.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.696"
$ws.Range("E50").Value = "  -11.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "29.91"
$ws.Range("E51").Value = "  -10.92%  "
